$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 handoff/handback datetimes advance to the next run
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 17:27:40"
$wsZhCn.Range("H2").Value = "2016-03-24 17:28:05"

# de-de sheet: row 2 handoff/handback datetimes advance to the next run
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 17:27:46"
$wsDeDe.Range("H2").Value = "2016-03-24 17:28:20"
